# Update the "Förändrad" date column (C) for rows 2-82 from 2026-02-22
# (serial 46075) to 2026-02-23 (serial 46076).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C82").Value = 46076
